# Re-pull / push data: update the dSF (column F) values to reflect the
# refreshed calculation for a handful of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -3
    5  = -8
    6  = 3
    8  = -6
    16 = 3
    17 = 3
    24 = -3
    29 = -10
    30 = -4
    35 = -3
    40 = -4
    42 = 0
    46 = 2
    47 = -4
    51 = 0
    53 = 5
    55 = 11
    60 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
